# Weekly update: a new price observation for
# "Vega Monumental Concepción - Berenjena" is inserted at row 94,
# pushing the previously-existing rows 94-97 down to rows 95-98.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 94 (shifts rows 94:97 -> 95:98).
$ws.Rows.Item(94).Insert()

# Populate the newly inserted row 94 with this week's data.
$ws.Cells.Item(94, 1).Value = 11
$ws.Cells.Item(94, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(94, 3).Value = "Bíobío"
$ws.Cells.Item(94, 4).Value = 44747
$ws.Cells.Item(94, 5).Value = 8
$ws.Cells.Item(94, 6).Value = 100112001
$ws.Cells.Item(94, 7).Value = "Berenjena"
$ws.Cells.Item(94, 8).Value = "Sin especificar"
$ws.Cells.Item(94, 9).Value = "Primera"
$ws.Cells.Item(94, 10).Value = 150
$ws.Cells.Item(94, 11).Value = 10000
$ws.Cells.Item(94, 12).Value = 12000
$ws.Cells.Item(94, 13).Value = 11067
$ws.Cells.Item(94, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(94, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(94, 16).Value = 184
$ws.Cells.Item(94, 17).Value = 60
$ws.Cells.Item(94, 18).Value = "Hortaliza"
